$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.683.26'
$ws.Range("E2").Value = '  +1.83%  '
$ws.Range("D3").Value = '2.212.11'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''267.75'
$ws.Range("E5").Value = '  +4.21%  '
$ws.Range("D6").Value = '''85.91'
$ws.Range("E6").Value = '  +10.92%  '
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").Value = '''46.03'
$ws.Range("E10").Value = '  +8.43%  '
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").Value = '''7.49'
$ws.Range("E12").Value = '  +6.15%  '
$ws.Range("D13").Value = '''0.104'
$ws.Range("E13").Value = '  +1.38%  '
$ws.Range("D14").Value = '2.543.49'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '''14.59'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = '2.207.03'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '''0.779'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '43.640.32'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '''5.97'
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").Value = '''69.82'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("E22").Value = '  +5.18%  '
$ws.Range("D23").Value = '''231.64'
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '''2.73'
$ws.Range("E24").Value = '  +23.95%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '''8.82'
$ws.Range("E25").Value = '  -5.87%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '''10.77'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("E28").Value = '  +5.65%  '
$ws.Range("D29").Value = '''39.10'
$ws.Range("E29").Value = '  -9.13%  '
$ws.Range("D30").Value = '''2.20'
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").Value = '''175.28'
$ws.Range("E31").Value = '  +1.33%  '
$ws.Range("D32").Value = '''0.0889'
$ws.Range("E32").Value = '  +1.49%  '
$ws.Range("D33").Value = '''20.52'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("E34").Value = '  +2.72%  '
$ws.Range("E35").Value = '  +1.60%  '
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("D37").Value = '''0.0356'
$ws.Range("E37").Value = '  -1.09%  '
$ws.Range("D38").Value = '''4.36'
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '''3.26'
$ws.Range("E39").Value = '  +15.08%  '
$ws.Range("D40").Value = '''12.29'
$ws.Range("E40").Value = '  -5.41%  '
$ws.Range("D41").Value = '''64.99'
$ws.Range("E41").Value = '  +7.73%  '
$ws.Range("E42").Value = '  -1.68%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.203'
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = '''5.40'
$ws.Range("E44").Value = '  +1.60%  '
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").Value = '''100.22'
$ws.Range("E46").Value = '  -2.93%  '
$ws.Range("D47").Value = '''8.32'
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").Value = '''0.440'
$ws.Range("E50").Value = '  -6.11%  '
$ws.Range("D51").Value = '''1.50'
$ws.Range("E51").Value = '  +3.77%  '
